$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Terry Rozier", "PG", "Miami Heat"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("RJ Barrett", "SG,SF,PF", "Toronto Raptors"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("T.J. McConnell", "PG", "Indiana Pacers"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Shai Gilgeous-Alexander", "PG,SG", "Oklahoma City Thunder"),
    @("Dennis Schröder", "PG", "Golden State Warriors")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row++
}
